$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.809499859809875
$ws.Range("B1").Value = 4.704745292663574
$ws.Range("C1").Value = 3.710319995880127
$ws.Range("D1").Value = 0.9018482565879822
$ws.Range("E1").Value = 0.4733693599700928
